$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for new log entries: Timestamp, Command, URL, Result, Entered Date, Entered Time
$rows = @(
    @('2024-09-23 18:48:34', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-23', '18:48:34'),
    @('2024-09-23 18:49:54', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-23', '18:49:54'),
    @('2024-09-23 18:50:19', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-23', '18:50:19'),
    @('2024-09-23 18:50:32', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:50:32'),
    @('2024-09-23 18:50:43', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-23', '18:50:43'),
    @('2024-09-23 18:50:50', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:50:50'),
    @('2024-09-23 18:51:07', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date current date is available for booking.', '2024-09-23', '18:51:07'),
    @('2024-09-23 18:51:10', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:51:10'),
    @('2024-09-23 18:56:06', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-wrapper button[aria-label*=''september 27'']"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:56:06'),
    @('2024-09-23 18:57:57', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:57:57'),
    @('2024-09-23 18:58:18', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:58:18'),
    @('2024-09-23 18:58:38', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '18:58:38'),
    @('2024-09-23 19:00:45', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '19:00:45'),
    @('2024-09-23 19:01:04', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '19:01:04'),
    @('2024-09-23 19:01:24', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '19:01:24'),
    @('2024-09-23 19:01:44', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.59); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF6705EFDA5+29557]
	(No symbol) [0x00007FF670562240]
	(No symbol) [0x00007FF67041B6EA]
	(No symbol) [0x00007FF67046FA15]
	(No symbol) [0x00007FF67046FC6C]
	(No symbol) [0x00007FF6704BBB07]
	(No symbol) [0x00007FF67049753F]
	(No symbol) [0x00007FF6704B88A3]
	(No symbol) [0x00007FF6704972A3]
	(No symbol) [0x00007FF6704612DF]
	(No symbol) [0x00007FF670462451]
	GetHandleVerifier [0x00007FF67091DCBD+3363469]
	GetHandleVerifier [0x00007FF670969B47+3674391]
	GetHandleVerifier [0x00007FF67095EAEB+3629243]
	GetHandleVerifier [0x00007FF6706AFC66+815670]
	(No symbol) [0x00007FF67056D6EF]
	(No symbol) [0x00007FF6705692B4]
	(No symbol) [0x00007FF670569450]
	(No symbol) [0x00007FF6705581FF]
	BaseThreadInitThunk [0x00007FFAEA46257D+29]
	RtlUserThreadStart [0x00007FFAEB4EAF28+40]
', '2024-09-23', '19:01:44'),
    @('2024-09-23 19:07:53', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date September 27 is available for booking.', '2024-09-23', '19:07:53'),
    @('2024-09-23 19:08:26', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Selected or default date September 27 is available for booking.', '2024-09-23', '19:08:26'),
    @('2024-09-23 19:10:04', 'check_availability', 'https://www.opentable.com/r/hals-the-steakhouse-nashville', 'Checked availability: Unable to determine availability. Please try again.', '2024-09-23', '19:10:04')
)

$startRow = 15
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]

    # Column E ("Entered Date") holds a YYYY-MM-DD string that must stay text,
    # not get auto-coerced into a date serial by Excel's input parser.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 5).Style = "Normal"

    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

